$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "Manguera4"
$ws.Range("I1").Value = "Elemento5"
[void]$ws.Range("L4").Select()
